$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B column accuracy values (rows 3-118) ---
$ws.Cells.Item(3, 2).Value = 0.96875
$ws.Cells.Item(4, 2).Value = 0.953125
$ws.Cells.Item(5, 2).Value = 0.9375
$ws.Cells.Item(6, 2).Value = 0.921875
$ws.Cells.Item(7, 2).Value = 0.90625
$ws.Cells.Item(8, 2).Value = 0.90625
$ws.Cells.Item(9, 2).Value = 0.90625
$ws.Cells.Item(10, 2).Value = 0.90625
$ws.Cells.Item(11, 2).Value = 0.90625
$ws.Cells.Item(12, 2).Value = 0.84375
$ws.Cells.Item(13, 2).Value = 0.875
$ws.Cells.Item(14, 2).Value = 0.859375
$ws.Cells.Item(15, 2).Value = 0.859375
$ws.Cells.Item(16, 2).Value = 0.84375
$ws.Cells.Item(17, 2).Value = 0.8125
$ws.Cells.Item(18, 2).Value = 0.84375
$ws.Cells.Item(19, 2).Value = 0.84375
$ws.Cells.Item(20, 2).Value = 0.796875
$ws.Cells.Item(21, 2).Value = 0.859375
$ws.Cells.Item(22, 2).Value = 0.796875
$ws.Cells.Item(23, 2).Value = 0.84375
$ws.Cells.Item(24, 2).Value = 0.84375
$ws.Cells.Item(25, 2).Value = 0.828125
$ws.Cells.Item(26, 2).Value = 0.828125
$ws.Cells.Item(27, 2).Value = 0.828125
$ws.Cells.Item(28, 2).Value = 0.828125
$ws.Cells.Item(29, 2).Value = 0.828125
$ws.Cells.Item(30, 2).Value = 0.828125
$ws.Cells.Item(31, 2).Value = 0.828125
$ws.Cells.Item(32, 2).Value = 0.828125
$ws.Cells.Item(33, 2).Value = 0.828125
$ws.Cells.Item(34, 2).Value = 0.828125
$ws.Cells.Item(35, 2).Value = 0.828125
$ws.Cells.Item(36, 2).Value = 0.828125
$ws.Cells.Item(37, 2).Value = 0.828125
$ws.Cells.Item(38, 2).Value = 0.828125
$ws.Cells.Item(39, 2).Value = 0.828125
$ws.Cells.Item(40, 2).Value = 0.828125
$ws.Cells.Item(41, 2).Value = 0.828125
$ws.Cells.Item(42, 2).Value = 0.828125
$ws.Cells.Item(43, 2).Value = 0.828125
$ws.Cells.Item(44, 2).Value = 0.828125
$ws.Cells.Item(45, 2).Value = 0.828125
$ws.Cells.Item(46, 2).Value = 0.828125
$ws.Cells.Item(47, 2).Value = 0.828125
$ws.Cells.Item(48, 2).Value = 0.828125
$ws.Cells.Item(49, 2).Value = 0.828125
$ws.Cells.Item(50, 2).Value = 0.828125
$ws.Cells.Item(51, 2).Value = 0.828125
$ws.Cells.Item(52, 2).Value = 0.796875
$ws.Cells.Item(53, 2).Value = 0.796875
$ws.Cells.Item(54, 2).Value = 0.796875
$ws.Cells.Item(55, 2).Value = 0.796875
$ws.Cells.Item(56, 2).Value = 0.796875
$ws.Cells.Item(57, 2).Value = 0.796875
$ws.Cells.Item(58, 2).Value = 0.796875
$ws.Cells.Item(59, 2).Value = 0.8125
$ws.Cells.Item(60, 2).Value = 0.796875
$ws.Cells.Item(61, 2).Value = 0.8125
$ws.Cells.Item(62, 2).Value = 0.8125
$ws.Cells.Item(63, 2).Value = 0.8125
$ws.Cells.Item(64, 2).Value = 0.828125
$ws.Cells.Item(65, 2).Value = 0.828125
$ws.Cells.Item(66, 2).Value = 0.828125
$ws.Cells.Item(67, 2).Value = 0.828125
$ws.Cells.Item(68, 2).Value = 0.828125
$ws.Cells.Item(69, 2).Value = 0.828125
$ws.Cells.Item(70, 2).Value = 0.828125
$ws.Cells.Item(71, 2).Value = 0.828125
$ws.Cells.Item(72, 2).Value = 0.828125
$ws.Cells.Item(73, 2).Value = 0.828125
$ws.Cells.Item(74, 2).Value = 0.828125
$ws.Cells.Item(75, 2).Value = 0.828125
$ws.Cells.Item(76, 2).Value = 0.828125
$ws.Cells.Item(77, 2).Value = 0.828125
$ws.Cells.Item(78, 2).Value = 0.828125
$ws.Cells.Item(79, 2).Value = 0.828125
$ws.Cells.Item(80, 2).Value = 0.828125
$ws.Cells.Item(81, 2).Value = 0.828125
$ws.Cells.Item(82, 2).Value = 0.828125
$ws.Cells.Item(83, 2).Value = 0.828125
$ws.Cells.Item(84, 2).Value = 0.828125
$ws.Cells.Item(85, 2).Value = 0.828125
$ws.Cells.Item(86, 2).Value = 0.828125
$ws.Cells.Item(87, 2).Value = 0.828125
$ws.Cells.Item(88, 2).Value = 0.828125
$ws.Cells.Item(89, 2).Value = 0.828125
$ws.Cells.Item(90, 2).Value = 0.828125
$ws.Cells.Item(91, 2).Value = 0.828125
$ws.Cells.Item(92, 2).Value = 0.828125
$ws.Cells.Item(93, 2).Value = 0.828125
$ws.Cells.Item(94, 2).Value = 0.828125
$ws.Cells.Item(95, 2).Value = 0.828125
$ws.Cells.Item(96, 2).Value = 0.828125
$ws.Cells.Item(97, 2).Value = 0.828125
$ws.Cells.Item(98, 2).Value = 0.828125
$ws.Cells.Item(99, 2).Value = 0.828125
$ws.Cells.Item(100, 2).Value = 0.828125
$ws.Cells.Item(101, 2).Value = 0.828125
$ws.Cells.Item(102, 2).Value = 0.828125
$ws.Cells.Item(103, 2).Value = 0.8125
$ws.Cells.Item(104, 2).Value = 0.796875
$ws.Cells.Item(105, 2).Value = 0.78125
$ws.Cells.Item(106, 2).Value = 0.796875
$ws.Cells.Item(107, 2).Value = 0.703125
$ws.Cells.Item(108, 2).Value = 0.796875
$ws.Cells.Item(109, 2).Value = 0.765625
$ws.Cells.Item(110, 2).Value = 0.765625
$ws.Cells.Item(111, 2).Value = 0.8125
$ws.Cells.Item(112, 2).Value = 0.8125
$ws.Cells.Item(113, 2).Value = 0.6875
$ws.Cells.Item(114, 2).Value = 0.765625
$ws.Cells.Item(115, 2).Value = 0.796875
$ws.Cells.Item(116, 2).Value = 0.796875
$ws.Cells.Item(117, 2).Value = 0.875
$ws.Cells.Item(118, 2).Value = 0.7377049180327869

# --- Update A column DisplayOutputs object repr text (rows 102-118) ---
$ws.Cells.Item(102, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(103, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(104, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(105, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(106, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(107, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(108, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(109, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(110, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(111, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(112, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(113, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(114, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(115, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(116, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(117, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'
$ws.Cells.Item(118, 1).Value = '<__main__.DisplayOutputs object at 0x7fece079ae20>'

# --- Update the sheet selection / active cell to match the author's last
#     interactive state (select A2:B118, active cell on B5) ---
$ws.Range("A2:B118").Select()
$ws.Range("B5").Activate()

$excel.CalculateFull()
